# Generate Report for handback
# - Marks the two source files as "Handed back" (was "Not yet handed off")
# - Records a handback: fills in "Latest Target File" / "Latest Handback File"
#   (columns E/F) with the same file references as the original source /
#   handoff file, and stamps a "Latest Handback DateTime" (column G).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$HandedBack = "Handed back"

# --- 1. Status column updates -------------------------------------------------

$overview.Range("B2").Value = $HandedBack
$overview.Range("C2").Value = $HandedBack
$overview.Range("B3").Value = $HandedBack
$overview.Range("C3").Value = $HandedBack

$zhcn.Range("B2").Value = $HandedBack
$zhcn.Range("B3").Value = $HandedBack

$dede.Range("B2").Value = $HandedBack
$dede.Range("B3").Value = $HandedBack

# --- helper: apply the same visual style used by the existing hyperlink cells

function Set-HyperlinkLook($range) {
    $range.Font.Underline = $True
    $range.Font.Color = 15570276
}

# --- 2. zh-cn sheet: Latest Target File (E) / Latest Handback File (F) -------

$mdUrl1040 = "https://github.com/OpenLocalizationTest/oltest/blob/4f19d056fbf0d701804e82f163bd187200857d19/e2e/1040d3c8-b286-49fe-b84b-360af021cc04.md"
$mdUrl629  = "https://github.com/OpenLocalizationTest/oltest/blob/4f19d056fbf0d701804e82f163bd187200857d19/e2e/629e121d-9986-4cbd-a3f7-576432f1d3cc.md"

$zhcnXlfUrl1040 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24a0ebddf5c06e6a3d457d2638a0a19ea11aa3f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/1040d3c8-b286-49fe-b84b-360af021cc04.42b36e6463acfeb9eaf13c93dba36098fb585622.zh-cn.xlf"
$zhcnXlfUrl629  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24a0ebddf5c06e6a3d457d2638a0a19ea11aa3f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/629e121d-9986-4cbd-a3f7-576432f1d3cc.88cdc26516c522a2ba95de41fd60c4412b68db3b.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("E2"), $mdUrl1040, [System.Type]::Missing, [System.Type]::Missing, "1040d3c8-b286-49fe-b84b-360af021cc04.md")
Set-HyperlinkLook $zhcn.Range("E2")

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhcnXlfUrl1040, [System.Type]::Missing, [System.Type]::Missing, "1040d3c8-b286-49fe-b84b-360af021cc04.42b36e6463acfeb9eaf13c93dba36098fb585622.zh-cn.xlf")
Set-HyperlinkLook $zhcn.Range("F2")

$zhcn.Hyperlinks.Add($zhcn.Range("E3"), $mdUrl629, [System.Type]::Missing, [System.Type]::Missing, "629e121d-9986-4cbd-a3f7-576432f1d3cc.md")
Set-HyperlinkLook $zhcn.Range("E3")

$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhcnXlfUrl629, [System.Type]::Missing, [System.Type]::Missing, "629e121d-9986-4cbd-a3f7-576432f1d3cc.88cdc26516c522a2ba95de41fd60c4412b68db3b.zh-cn.xlf")
Set-HyperlinkLook $zhcn.Range("F3")

# Latest Handback DateTime (column G)
$zhcn.Range("G2").Value = "2016-01-08 15:24:44"
$zhcn.Range("G3").Value = "2016-01-08 15:24:44"

# --- 3. de-de sheet: Latest Target File (E) / Latest Handback File (F) ------

$dedeXlfUrl1040 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b895b1175b9edb9bbb6ec945a3d1f96c77823c31/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/1040d3c8-b286-49fe-b84b-360af021cc04.42b36e6463acfeb9eaf13c93dba36098fb585622.de-de.xlf"
$dedeXlfUrl629  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b895b1175b9edb9bbb6ec945a3d1f96c77823c31/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/629e121d-9986-4cbd-a3f7-576432f1d3cc.88cdc26516c522a2ba95de41fd60c4412b68db3b.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("E2"), $mdUrl1040, [System.Type]::Missing, [System.Type]::Missing, "1040d3c8-b286-49fe-b84b-360af021cc04.md")
Set-HyperlinkLook $dede.Range("E2")

$dede.Hyperlinks.Add($dede.Range("F2"), $dedeXlfUrl1040, [System.Type]::Missing, [System.Type]::Missing, "1040d3c8-b286-49fe-b84b-360af021cc04.42b36e6463acfeb9eaf13c93dba36098fb585622.de-de.xlf")
Set-HyperlinkLook $dede.Range("F2")

$dede.Hyperlinks.Add($dede.Range("E3"), $mdUrl629, [System.Type]::Missing, [System.Type]::Missing, "629e121d-9986-4cbd-a3f7-576432f1d3cc.md")
Set-HyperlinkLook $dede.Range("E3")

$dede.Hyperlinks.Add($dede.Range("F3"), $dedeXlfUrl629, [System.Type]::Missing, [System.Type]::Missing, "629e121d-9986-4cbd-a3f7-576432f1d3cc.88cdc26516c522a2ba95de41fd60c4412b68db3b.de-de.xlf")
Set-HyperlinkLook $dede.Range("F3")

# Latest Handback DateTime (column G)
$dede.Range("G2").Value = "2016-01-08 15:25:08"
$dede.Range("G3").Value = "2016-01-08 15:25:08"

Write-Output "handback report generated"
